# Updated cryptos list (price/volume refresh) matching the GitHub Actions commit.
# Numeric-looking "Price" strings are entered with a leading apostrophe so Excel
# keeps them as text (matching the original inlineStr cells) instead of coercing
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.548.19'
$ws.Range('E2').Value = '  -0.50%  '

$ws.Range('D3').Value = '1.623.46'
$ws.Range('E3').Value = '  -1.25%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = "'211.56"

$ws.Range('D6').Value = "'0.525"
$ws.Range('E6').Value = '  -0.82%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = "'23.19"
$ws.Range('E8').Value = '  -0.45%  '

$ws.Range('E9').Value = '  +1.53%  '

$ws.Range('E10').Value = '  -0.17%  '

$ws.Range('D11').Value = "'0.0884"
$ws.Range('E11').Value = '  -1.06%  '

$ws.Range('D12').Value = '1.854.62'
$ws.Range('E12').Value = '  -1.18%  '

$ws.Range('D13').Value = '1.626.21'
$ws.Range('E13').Value = '  -1.18%  '

$ws.Range('E14').Value = '  -0.03%  '

$ws.Range('E15').Value = '  -2.18%  '

$ws.Range('D16').Value = "'65.07"
$ws.Range('E16').Value = '  +0.70%  '

$ws.Range('D17').Value = '27.512.88'
$ws.Range('E17').Value = '  -0.53%  '

$ws.Range('D18').Value = "'230.42"
$ws.Range('E18').Value = '  -0.14%  '

$ws.Range('D19').Value = '0.0₃0719'
$ws.Range('E19').Value = '  -0.61%  '

$ws.Range('D20').Value = "'7.52"
$ws.Range('E20').Value = '  -1.41%  '

$ws.Range('E21').Value = '  -0.05%  '

$ws.Range('D22').Value = "'10.34"
$ws.Range('E22').Value = '  +3.38%  '

$ws.Range('E23').Value = '  +0.93%  '

$ws.Range('E24').Value = '  +5.77%  '

$ws.Range('D25').Value = "'148.71"
$ws.Range('E25').Value = '  -0.80%  '

$ws.Range('D26').Value = "'6.87"
$ws.Range('E26').Value = '  -0.90%  '

$ws.Range('E27').Value = '  -0.80%  '

$ws.Range('E29').Value = '  -0.70%  '

$ws.Range('D30').Value = "'1.17"
$ws.Range('E30').Value = '  -0.99%  '

$ws.Range('D31').Value = "'0.0483"
$ws.Range('E31').Value = '  -0.88%  '

$ws.Range('E32').Value = '  -0.88%  '

$ws.Range('D33').Value = '1.467.03'
$ws.Range('E33').Value = '  +1.62%  '

$ws.Range('D34').Value = "'3.05"
$ws.Range('E34').Value = '  -2.82%  '

$ws.Range('E35').Value = '  -3.05%  '

$ws.Range('D36').Value = "'2.34"
$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('D37').Value = "'0.932"
$ws.Range('E37').Value = '  +5.21%  '

$ws.Range('D38').Value = "'0.873"
$ws.Range('E38').Value = '  -0.79%  '

$ws.Range('D39').Value = "'0.555"
$ws.Range('E39').Value = '  -2.29%  '

$ws.Range('E40').Value = '  -0.25%  '

$ws.Range('E41').Value = '  +0.00%  '

$ws.Range('E42').Value = '  -2.12%  '

$ws.Range('D43').Value = "'67.24"
$ws.Range('E43').Value = '  +0.18%  '

$ws.Range('E44').Value = '  -1.61%  '

$ws.Range('E45').Value = '  -1.62%  '

$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'5.28"
$ws.Range('E46').Value = '  -6.14%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'1.76"
$ws.Range('E47').Value = '  +1.92%  '

$ws.Range('D48').Value = '1.763.83'
$ws.Range('E48').Value = '  -1.29%  '

$ws.Range('D49').Value = "'87.35"
$ws.Range('E49').Value = '  +2.02%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0106'
$ws.Range('E50').Value = '  -1.22%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.0999"
$ws.Range('E51').Value = '  +1.09%  '
